# "merged with latest from UW team's update"
#
# This updates the XLSForm "selects" workbook:
#  - survey sheet: the remote-data cascading select now starts with a
#    "countries" level (feeding into the existing "states" level, replacing
#    the old "cities" level); a new "Which devices do you use?" note plus a
#    4-choice (desktop/laptop/smartphone/tablet) inline select_one replaces
#    the old 3-choice (i1/i2/i3) placeholder select; and the content-provider
#    question gets an explicit (false) relevant condition plus two trailing
#    blank rows.
#  - queries sheet: the "states"/"cities" callback rows are renamed to
#    "countries"/"states" and the states query now filters on the country.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# survey sheet
# ---------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

# Insert the new "Which devices do you use?" note row above the old inline
# select block (old row 10), and a 4th row inside that block (old row 13,
# now row 14) so the 3-choice select becomes a 4-choice select.
$survey.Rows.Item(10).Insert()
$survey.Rows.Item(14).Insert()

# Row 2: remote-data cascading select now begins with "country" instead of
# "state" (the old "state" level moves down to row 3).
$survey.Range("B2").Value = "select_one countries"
$survey.Range("F2").Value = "country"
$survey.Range("G2").Value = "Choose a country:"

# Row 3: what used to be the "city" level is now the "state" level, reusing
# the same query name ("states") as before.
$survey.Range("B3").Value = "select_one states"
$survey.Range("F3").Value = "state"
$survey.Range("G3").Value = "Choose a state:"

# Row 10 (new): plain note introducing the device question.
$survey.Range("B10").Value = "note"
$survey.Range("G10").Value = "Which devices do you use?"

# Rows 11-14: inline select_one yes_no -> inline select_one with four device
# choices.
$survey.Range("B11").Value = "select_one yes_no"
$survey.Range("C11").Value = "inline"
$survey.Range("F11").Value = "desktop"
$survey.Range("G11").Value = "Desktop computer"

$survey.Range("B12").Value = "select_one yes_no"
$survey.Range("C12").Value = "inline"
$survey.Range("F12").Value = "laptop"
$survey.Range("G12").Value = "Laptop computer"

$survey.Range("B13").Value = "select_one yes_no"
$survey.Range("C13").Value = "inline"
$survey.Range("F13").Value = "smartphone"
$survey.Range("G13").Value = "Smartphone"

$survey.Range("B14").Value = "select_one yes_no"
$survey.Range("C14").Value = "inline"
$survey.Range("F14").Value = "tablet"
$survey.Range("G14").Value = "Tablet"

# Row 21 (content provider test, shifted down by the two inserted rows)
# gains an explicit FALSE relevant condition.
$survey.Range("E21").Value = $false

# Column E got a bit narrower.
$survey.Columns.Item(5).ColumnWidth = 33.29

# ---------------------------------------------------------------------
# queries sheet
# ---------------------------------------------------------------------
$queries = $wb.Worksheets.Item("queries")

$queries.Range("A2").Value = "countries"
$queries.Range("A3").Value = "states"
$queries.Range("B3").Value = '"https://query.yahooapis.com/v1/public/yql?format=json&q=" +  encodeURIComponent("select * from geo.states where place=''" + data(''country'') + "''")'
$queries.Range("C6").Value = "context"
